$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 23-24, shifting existing rows 23-58 down to 25-60
$ws.Rows("23:24").Insert()

# Populate new row 23
$ws.Cells.Item(23,1).Value = 5
$ws.Cells.Item(23,2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(23,3).Value = 'Maule'
$ws.Cells.Item(23,4).Value = 45280
$ws.Cells.Item(23,5).Value = 7
$ws.Cells.Item(23,6).Value = 'Fruta'
$ws.Cells.Item(23,7).Value = 100103
$ws.Cells.Item(23,8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(23,9).Value = 100103003
$ws.Cells.Item(23,10).Value = 'Damasco'
$ws.Cells.Item(23,11).Value = 'Castle Brite'
$ws.Cells.Item(23,12).Value = 'Primera'
$ws.Cells.Item(23,13).Value = 100
$ws.Cells.Item(23,14).Value = 13000
$ws.Cells.Item(23,15).Value = 13000
$ws.Cells.Item(23,16).Value = 13000
$ws.Cells.Item(23,17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(23,18).Value = 'Región de O''Higgins'
$ws.Cells.Item(23,19).Value = 1300
$ws.Cells.Item(23,20).Value = 10

# Populate new row 24
$ws.Cells.Item(24,1).Value = 5
$ws.Cells.Item(24,2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(24,3).Value = 'Maule'
$ws.Cells.Item(24,4).Value = 45280
$ws.Cells.Item(24,5).Value = 7
$ws.Cells.Item(24,6).Value = 'Fruta'
$ws.Cells.Item(24,7).Value = 100103
$ws.Cells.Item(24,8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(24,9).Value = 100103003
$ws.Cells.Item(24,10).Value = 'Damasco'
$ws.Cells.Item(24,11).Value = 'Castle Brite'
$ws.Cells.Item(24,12).Value = 'Primera'
$ws.Cells.Item(24,13).Value = 180
$ws.Cells.Item(24,14).Value = 16000
$ws.Cells.Item(24,15).Value = 16000
$ws.Cells.Item(24,16).Value = 16000
$ws.Cells.Item(24,17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(24,18).Value = 'Provincia de Limarí'
$ws.Cells.Item(24,19).Value = 1000
$ws.Cells.Item(24,20).Value = 16
